# Auto-generated Excel COM-interop script to update Typhon Profits leve data
# Applies scheduled market-price refresh values to the H:N columns across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70: columns H,I,J,K,L,M,N
$ws.Cells.Item(70, 8).Value = 1195
$ws.Cells.Item(70, 9).Value = 1390
$ws.Cells.Item(70, 10).Value = 1000
$ws.Cells.Item(70, 11).Value = 4170
$ws.Cells.Item(70, 12).Value = 3000
$ws.Cells.Item(70, 13).Value = -3900
$ws.Cells.Item(70, 14).Value = -3540
# Row 73: columns H,I,J,K,L,M,N
$ws.Cells.Item(73, 8).Value = 1195
$ws.Cells.Item(73, 9).Value = 1390
$ws.Cells.Item(73, 10).Value = 1000
$ws.Cells.Item(73, 11).Value = 4170
$ws.Cells.Item(73, 12).Value = 3000
$ws.Cells.Item(73, 13).Value = -3234
$ws.Cells.Item(73, 14).Value = -4872
# Row 95: columns H,J,L,N
$ws.Cells.Item(95, 8).Value = 27874.666
$ws.Cells.Item(95, 10).Value = 27874.666
$ws.Cells.Item(95, 12).Value = 27874.666
$ws.Cells.Item(95, 14).Value = -33366.666
# Row 99: columns H,I,J,K,L,M,N
$ws.Cells.Item(99, 8).Value = 233.33333
$ws.Cells.Item(99, 9).Value = 200
$ws.Cells.Item(99, 10).Value = 300
$ws.Cells.Item(99, 11).Value = 600
$ws.Cells.Item(99, 12).Value = 900
$ws.Cells.Item(99, 13).Value = 898
$ws.Cells.Item(99, 14).Value = -3896
# Row 100: columns H,I,J,K,L,M,N
$ws.Cells.Item(100, 8).Value = 2182.4
$ws.Cells.Item(100, 9).Value = 1604.3636
$ws.Cells.Item(100, 10).Value = 2888.889
$ws.Cells.Item(100, 11).Value = 1604.3636
$ws.Cells.Item(100, 12).Value = 2888.889
$ws.Cells.Item(100, 13).Value = -1063.3636
$ws.Cells.Item(100, 14).Value = -3970.889
# Row 116: columns H,I,J,K,L,M,N
$ws.Cells.Item(116, 8).Value = 18821368
$ws.Cells.Item(116, 9).Value = 47044370
$ws.Cells.Item(116, 10).Value = 6034
$ws.Cells.Item(116, 11).Value = 47044370
$ws.Cells.Item(116, 12).Value = 6034
$ws.Cells.Item(116, 13).Value = -47040928
$ws.Cells.Item(116, 14).Value = -12918
# Row 138: columns H,J,L,N
$ws.Cells.Item(138, 8).Value = 3882.568
$ws.Cells.Item(138, 10).Value = 3864.0264
$ws.Cells.Item(138, 12).Value = 11592.0792
$ws.Cells.Item(138, 14).Value = -21872.0792

$ws = $wb.Worksheets.Item("ARM")
# Row 2: columns H,I,K,M
$ws.Cells.Item(2, 8).Value = 1074.5
$ws.Cells.Item(2, 9).Value = 1031.7059
$ws.Cells.Item(2, 11).Value = 1031.7059
$ws.Cells.Item(2, 13).Value = -918.7058999999999
# Row 45: columns H,I,J,K,L,M,N
$ws.Cells.Item(45, 8).Value = 2405.7837
$ws.Cells.Item(45, 9).Value = 2092.074
$ws.Cells.Item(45, 10).Value = 3252.8
$ws.Cells.Item(45, 11).Value = 2092.074
$ws.Cells.Item(45, 12).Value = 3252.8
$ws.Cells.Item(45, 13).Value = -1715.074
$ws.Cells.Item(45, 14).Value = -4006.8
# Row 116: columns H,I,K,M
$ws.Cells.Item(116, 8).Value = 1074.5
$ws.Cells.Item(116, 9).Value = 1031.7059
$ws.Cells.Item(116, 11).Value = 1031.7059
$ws.Cells.Item(116, 13).Value = 1262.2941
# Row 122: columns H,I,K,M
$ws.Cells.Item(122, 8).Value = 1872.2927
$ws.Cells.Item(122, 9).Value = 1768.25
$ws.Cells.Item(122, 11).Value = 5304.75
$ws.Cells.Item(122, 13).Value = -2854.75

$ws = $wb.Worksheets.Item("BSM")
# Row 3: columns H,I,K,M
$ws.Cells.Item(3, 8).Value = 1074.5
$ws.Cells.Item(3, 9).Value = 1031.7059
$ws.Cells.Item(3, 11).Value = 1031.7059
$ws.Cells.Item(3, 13).Value = -917.7058999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31: columns H,I,J,K,L,M,N
$ws.Cells.Item(31, 8).Value = 3787.8447
$ws.Cells.Item(31, 9).Value = 1788.1765
$ws.Cells.Item(31, 10).Value = 6620.7085
$ws.Cells.Item(31, 11).Value = 1788.1765
$ws.Cells.Item(31, 12).Value = 6620.7085
$ws.Cells.Item(31, 13).Value = -1493.1765
$ws.Cells.Item(31, 14).Value = -7210.7085
# Row 34: columns H,I,J,K,L,M,N
$ws.Cells.Item(34, 8).Value = 3787.8447
$ws.Cells.Item(34, 9).Value = 1788.1765
$ws.Cells.Item(34, 10).Value = 6620.7085
$ws.Cells.Item(34, 11).Value = 1788.1765
$ws.Cells.Item(34, 12).Value = 6620.7085
$ws.Cells.Item(34, 13).Value = -1586.1765
$ws.Cells.Item(34, 14).Value = -7024.7085
# Row 74: columns H,J,L,N
$ws.Cells.Item(74, 8).Value = 42622.9
$ws.Cells.Item(74, 10).Value = 42622.9
$ws.Cells.Item(74, 12).Value = 42622.9
$ws.Cells.Item(74, 14).Value = -44370.9
# Row 77: columns H,J,L,N
$ws.Cells.Item(77, 8).Value = 42622.9
$ws.Cells.Item(77, 10).Value = 42622.9
$ws.Cells.Item(77, 12).Value = 127868.7
$ws.Cells.Item(77, 14).Value = -136604.7
# Row 99: columns H,I,J,K,L,M,N
$ws.Cells.Item(99, 8).Value = 3456.1714
$ws.Cells.Item(99, 9).Value = 2913.889
$ws.Cells.Item(99, 10).Value = 4030.353
$ws.Cells.Item(99, 11).Value = 2913.889
$ws.Cells.Item(99, 12).Value = 4030.353
$ws.Cells.Item(99, 13).Value = -1415.889
$ws.Cells.Item(99, 14).Value = -7026.353
# Row 122: columns H,I,J,K,L,M,N
$ws.Cells.Item(122, 8).Value = 1014.7778
$ws.Cells.Item(122, 9).Value = 849.25
$ws.Cells.Item(122, 10).Value = 1255.5454
$ws.Cells.Item(122, 11).Value = 2547.75
$ws.Cells.Item(122, 12).Value = 3766.6362
$ws.Cells.Item(122, 13).Value = -97.75
$ws.Cells.Item(122, 14).Value = -8666.636200000001
# Row 126: columns H,I,J,K,L,M,N
$ws.Cells.Item(126, 8).Value = 3456.1714
$ws.Cells.Item(126, 9).Value = 2913.889
$ws.Cells.Item(126, 10).Value = 4030.353
$ws.Cells.Item(126, 11).Value = 8741.667000000001
$ws.Cells.Item(126, 12).Value = 12091.059
$ws.Cells.Item(126, 13).Value = -6271.667000000001
$ws.Cells.Item(126, 14).Value = -17031.059
# Row 132: columns H,I,J,K,L,M,N
$ws.Cells.Item(132, 8).Value = 50003460
$ws.Cells.Item(132, 9).Value = 76925624
$ws.Cells.Item(132, 10).Value = 5149.7144
$ws.Cells.Item(132, 11).Value = 230776872
$ws.Cells.Item(132, 12).Value = 15449.1432
$ws.Cells.Item(132, 13).Value = -230774342
$ws.Cells.Item(132, 14).Value = -20509.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 5: columns H,I,K,M
$ws.Cells.Item(5, 8).Value = 1271.6487
$ws.Cells.Item(5, 9).Value = 834.7308
$ws.Cells.Item(5, 11).Value = 2504.1924
$ws.Cells.Item(5, 13).Value = -2392.1924
# Row 129: columns H,J,L,N
$ws.Cells.Item(129, 8).Value = 224252.22
$ws.Cells.Item(129, 10).Value = 268992
$ws.Cells.Item(129, 12).Value = 806976
$ws.Cells.Item(129, 14).Value = -816976
# Row 131: columns H,I,J,K,L,M,N
$ws.Cells.Item(131, 8).Value = 838.1622
$ws.Cells.Item(131, 9).Value = 476.5
$ws.Cells.Item(131, 10).Value = 894.6719000000001
$ws.Cells.Item(131, 11).Value = 1429.5
$ws.Cells.Item(131, 12).Value = 2684.0157
$ws.Cells.Item(131, 13).Value = 3610.5
$ws.Cells.Item(131, 14).Value = -12764.0157
# Row 135: columns H,I,K,M
$ws.Cells.Item(135, 8).Value = 1271.6487
$ws.Cells.Item(135, 9).Value = 834.7308
$ws.Cells.Item(135, 11).Value = 7512.577200000001
$ws.Cells.Item(135, 13).Value = -4977.577200000001

$ws = $wb.Worksheets.Item("GSM")
# Row 95: columns H,J,L,N
$ws.Cells.Item(95, 8).Value = 25000
$ws.Cells.Item(95, 10).Value = 25000
$ws.Cells.Item(95, 12).Value = 25000
$ws.Cells.Item(95, 14).Value = -30492
# Row 102: columns H,I,K,M
$ws.Cells.Item(102, 8).Value = 3040.647
$ws.Cells.Item(102, 9).Value = 2710.8667
$ws.Cells.Item(102, 11).Value = 2710.8667
$ws.Cells.Item(102, 13).Value = -1088.8667
# Row 122: columns H,I,J,K,L,M,N
$ws.Cells.Item(122, 8).Value = 6592.143
$ws.Cells.Item(122, 9).Value = 7786.25
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 23358.75
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -20908.75
$ws.Cells.Item(122, 14).Value = -19900

$ws = $wb.Worksheets.Item("LTW")
# Row 7: columns H,I,J,K,L,M,N
$ws.Cells.Item(7, 8).Value = 3443.75
$ws.Cells.Item(7, 9).Value = 3108.3333
$ws.Cells.Item(7, 10).Value = 4450
$ws.Cells.Item(7, 11).Value = 3108.3333
$ws.Cells.Item(7, 12).Value = 4450
$ws.Cells.Item(7, 13).Value = -2996.3333
$ws.Cells.Item(7, 14).Value = -4674
# Row 40: columns H,I,K,M
$ws.Cells.Item(40, 8).Value = 4494.591
$ws.Cells.Item(40, 9).Value = 4280.048
$ws.Cells.Item(40, 11).Value = 4280.048
$ws.Cells.Item(40, 13).Value = -4144.048
# Row 61: columns H,I,J,K,L,M,N
$ws.Cells.Item(61, 8).Value = 4069.5881
$ws.Cells.Item(61, 9).Value = 2197.923
$ws.Cells.Item(61, 10).Value = 10152.5
$ws.Cells.Item(61, 11).Value = 2197.923
$ws.Cells.Item(61, 12).Value = 10152.5
$ws.Cells.Item(61, 13).Value = -1995.923
$ws.Cells.Item(61, 14).Value = -10556.5
# Row 68: columns H,I,J,K,L,M,N
$ws.Cells.Item(68, 8).Value = 2770.7144
$ws.Cells.Item(68, 9).Value = 2399.5
$ws.Cells.Item(68, 10).Value = 2919.2
$ws.Cells.Item(68, 11).Value = 2399.5
$ws.Cells.Item(68, 12).Value = 2919.2
$ws.Cells.Item(68, 13).Value = -1650.5
$ws.Cells.Item(68, 14).Value = -4417.2
# Row 71: columns H,I,J,K,L,M,N
$ws.Cells.Item(71, 8).Value = 2770.7144
$ws.Cells.Item(71, 9).Value = 2399.5
$ws.Cells.Item(71, 10).Value = 2919.2
$ws.Cells.Item(71, 11).Value = 11997.5
$ws.Cells.Item(71, 12).Value = 14596
$ws.Cells.Item(71, 13).Value = -8253.5
$ws.Cells.Item(71, 14).Value = -22084
# Row 113: columns H,I,J,K,L,M,N
$ws.Cells.Item(113, 8).Value = 4069.5881
$ws.Cells.Item(113, 9).Value = 2197.923
$ws.Cells.Item(113, 10).Value = 10152.5
$ws.Cells.Item(113, 11).Value = 2197.923
$ws.Cells.Item(113, 12).Value = 10152.5
$ws.Cells.Item(113, 13).Value = -27.92299999999977
$ws.Cells.Item(113, 14).Value = -14492.5
# Row 126: columns H,I,J,K,L,M,N
$ws.Cells.Item(126, 8).Value = 3443.75
$ws.Cells.Item(126, 9).Value = 3108.3333
$ws.Cells.Item(126, 10).Value = 4450
$ws.Cells.Item(126, 11).Value = 9324.999899999999
$ws.Cells.Item(126, 12).Value = 13350
$ws.Cells.Item(126, 13).Value = -6854.999899999999
$ws.Cells.Item(126, 14).Value = -18290
# Row 136: columns H,I,J,K,L,M,N
$ws.Cells.Item(136, 8).Value = 2191.64
$ws.Cells.Item(136, 9).Value = 2186.5652
$ws.Cells.Item(136, 10).Value = 2250
$ws.Cells.Item(136, 11).Value = 6559.6956
$ws.Cells.Item(136, 12).Value = 6750
$ws.Cells.Item(136, 13).Value = -4009.6956
$ws.Cells.Item(136, 14).Value = -11850

$ws = $wb.Worksheets.Item("WVR")
# Row 62: columns H,I,J,K,L,M,N
$ws.Cells.Item(62, 8).Value = 4700
$ws.Cells.Item(62, 9).Value = 4000
$ws.Cells.Item(62, 10).Value = 4980
$ws.Cells.Item(62, 11).Value = 4000
$ws.Cells.Item(62, 12).Value = 4980
$ws.Cells.Item(62, 13).Value = -3376
$ws.Cells.Item(62, 14).Value = -6228
# Row 65: columns H,I,J,K,L,M,N
$ws.Cells.Item(65, 8).Value = 4700
$ws.Cells.Item(65, 9).Value = 4000
$ws.Cells.Item(65, 10).Value = 4980
$ws.Cells.Item(65, 11).Value = 20000
$ws.Cells.Item(65, 12).Value = 24900
$ws.Cells.Item(65, 13).Value = -16880
$ws.Cells.Item(65, 14).Value = -31140
# Row 76: columns H,J,L,N
$ws.Cells.Item(76, 8).Value = 32000
$ws.Cells.Item(76, 10).Value = 32000
$ws.Cells.Item(76, 12).Value = 32000
$ws.Cells.Item(76, 14).Value = -32630
# Row 79: columns H,J,L,N
$ws.Cells.Item(79, 8).Value = 32000
$ws.Cells.Item(79, 10).Value = 32000
$ws.Cells.Item(79, 12).Value = 32000
$ws.Cells.Item(79, 14).Value = -34184
# Row 95: columns H,J,L,N
$ws.Cells.Item(95, 8).Value = 20172
$ws.Cells.Item(95, 10).Value = 20172
$ws.Cells.Item(95, 12).Value = 20172
$ws.Cells.Item(95, 14).Value = -25664
# Row 122: columns H,I,K,M
$ws.Cells.Item(122, 8).Value = 1883.3334
$ws.Cells.Item(122, 9).Value = 1800
$ws.Cells.Item(122, 11).Value = 5400
$ws.Cells.Item(122, 13).Value = -2950

